# Auto-generated Excel COM-interop edit script
# Applies the 2026-02-10 17:20 automatic data refresh to Dades_Meteo sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-10 17:18:41"
$ws.Range("I2").Value = "29.4 mm"
$ws.Range("O2").Value = "0.7 °C"
$ws.Range("E3").Value = "2026-02-10 17:18:43"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "96%"
$ws.Range("H3").NumberFormat = "General"
$ws.Range("I3").Value = "18.4 mm"
$ws.Range("L3").Value = "47.5 km/h - 229º 16:54 TU"
$ws.Range("O3").Value = "0.8 °C"
$ws.Range("E4").Value = "2026-02-10 17:18:46"
$ws.Range("J4").Value = "1004.2 hPa"
$ws.Range("O4").Value = "11.7 °C"
$ws.Range("E5").Value = "2026-02-10 17:18:48"
$ws.Range("I5").Value = "24.7 mm"
$ws.Range("O5").Value = "1.4 °C"
$ws.Range("E6").Value = "2026-02-10 17:18:50"
$ws.Range("J6").Value = "1004.5 hPa"
$ws.Range("L6").Value = "13.0 km/h - 38º 16:41 TU"
$ws.Range("O6").Value = "9.5 °C"
$ws.Range("E7").Value = "2026-02-10 17:18:53"
$ws.Range("J7").Value = "1004.9 hPa"
$ws.Range("O7").Value = "14.8 °C"
$ws.Range("E8").Value = "2026-02-10 17:18:55"
$ws.Range("J8").Value = "1004.8 hPa"
$ws.Range("O8").Value = "11.5 °C"
$ws.Range("E9").Value = "2026-02-10 17:18:58"
$ws.Range("O9").Value = "8.7 °C"
$ws.Range("E10").Value = "2026-02-10 17:19:00"
$ws.Range("O10").Value = "9.9 °C"
$ws.Range("E11").Value = "2026-02-10 17:19:02"
$ws.Range("O11").Value = "7.2 °C"
$ws.Range("E12").Value = "2026-02-10 17:19:04"
$ws.Range("O12").Value = "8.8 °C"
$ws.Range("E13").Value = "2026-02-10 17:19:07"
$ws.Range("J13").Value = "1007.2 hPa"
$ws.Range("K13").Value = "7.6 MJ/m2"
$ws.Range("O13").Value = "4.8 °C"
$ws.Range("E14").Value = "2026-02-10 17:19:10"
$ws.Range("O14").Value = "13.1 °C"
$ws.Range("E15").Value = "2026-02-10 17:19:12"
$ws.Range("E16").Value = "2026-02-10 17:19:14"
$ws.Range("I16").Value = "20.1 mm"
$ws.Range("E17").Value = "2026-02-10 17:19:16"
$ws.Range("E18").Value = "2026-02-10 17:19:19"
$ws.Range("O18").Value = "10.0 °C"
$ws.Range("E19").Value = "2026-02-10 17:19:22"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "94%"
$ws.Range("H19").NumberFormat = "General"
$ws.Range("K19").Value = "5.9 MJ/m2"
$ws.Range("L19").Value = "19.4 km/h - 264º 16:48 TU"
$ws.Range("O19").Value = "6.1 °C"
$ws.Range("E20").Value = "2026-02-10 17:19:24"
$ws.Range("I20").Value = "3.3 mm"
$ws.Range("K20").Value = "7.1 MJ/m2"
$ws.Range("M20").Value = "2.4 °C 16:52 TU"
$ws.Range("O20").Value = "0.1 °C"
$ws.Range("E21").Value = "2026-02-10 17:19:27"
$ws.Range("J21").Value = "1006.5 hPa"
$ws.Range("O21").Value = "6.9 °C"
$ws.Range("E22").Value = "2026-02-10 17:19:29"
$ws.Range("I22").Value = "7.0 mm"
$ws.Range("O22").Value = "-0.7 °C"
$ws.Range("E23").Value = "2026-02-10 17:19:32"
$ws.Range("E24").Value = "2026-02-10 17:19:34"
$ws.Range("J24").Value = "1006.4 hPa"
$ws.Range("K24").Value = "9.9 MJ/m2"
$ws.Range("O24").Value = "11.0 °C"
$ws.Range("E25").Value = "2026-02-10 17:19:37"
$ws.Range("K25").Value = "6.9 MJ/m2"
$ws.Range("E26").Value = "2026-02-10 17:19:40"
$ws.Range("J26").Value = "1003.6 hPa"
$ws.Range("O26").Value = "6.0 °C"
$ws.Range("E27").Value = "2026-02-10 17:19:42"
$ws.Range("K27").Value = "7.1 MJ/m2"
$ws.Range("E28").Value = "2026-02-10 17:19:45"
$ws.Range("J28").Value = "1004.6 hPa"
$ws.Range("L28").Value = "15.1 km/h - 70º 16:47 TU"
$ws.Range("O28").Value = "8.6 °C"
$ws.Range("E29").Value = "2026-02-10 17:19:48"
$ws.Range("O29").Value = "10.4 °C"
$ws.Range("E30").Value = "2026-02-10 17:19:50"
$ws.Range("J30").Value = "1004.6 hPa"
$ws.Range("O30").Value = "9.1 °C"
$ws.Range("E31").Value = "2026-02-10 17:19:53"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "80%"
$ws.Range("H31").NumberFormat = "General"
$ws.Range("I31").Value = "0.7 mm"
$ws.Range("J31").Value = "1003.8 hPa"
$ws.Range("E32").Value = "2026-02-10 17:19:56"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "92%"
$ws.Range("H32").NumberFormat = "General"
$ws.Range("K32").Value = "8.2 MJ/m2"
$ws.Range("E33").Value = "2026-02-10 17:19:59"
$ws.Range("J33").Value = "1006.9 hPa"
$ws.Range("O33").Value = "3.8 °C"
$ws.Range("E34").Value = "2026-02-10 17:20:01"
$ws.Range("O34").Value = "3.7 °C"
$ws.Range("E35").Value = "2026-02-10 17:20:04"
$ws.Range("K35").Value = "10.5 MJ/m2"
$ws.Range("O35").Value = "12.8 °C"
$ws.Range("E36").Value = "2026-02-10 17:20:07"
$ws.Range("I36").Value = "0.4 mm"
$ws.Range("J36").Value = "1004.7 hPa"
$ws.Range("E37").Value = "2026-02-10 17:20:09"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "90%"
$ws.Range("H37").NumberFormat = "General"
$ws.Range("J37").Value = "1006.0 hPa"
$ws.Range("O37").Value = "6.3 °C"
$ws.Range("E38").Value = "2026-02-10 17:20:12"
$ws.Range("K38").Value = "8.0 MJ/m2"
$ws.Range("O38").Value = "10.4 °C"
$ws.Range("E39").Value = "2026-02-10 17:20:14"
$ws.Range("K39").Value = "8.1 MJ/m2"
$ws.Range("M39").Value = "4.5 °C 16:56 TU"
$ws.Range("E40").Value = "2026-02-10 17:20:17"
$ws.Range("I40").Value = "5.5 mm"
$ws.Range("J40").Value = "1007.4 hPa"
$ws.Range("O40").Value = "7.3 °C"
$ws.Range("E41").Value = "2026-02-10 17:20:20"
$ws.Range("J41").Value = "1005.0 hPa"
$ws.Range("K41").Value = "9.6 MJ/m2"
$ws.Range("O41").Value = "14.1 °C"
$ws.Range("E42").Value = "2026-02-10 17:20:22"
$ws.Range("O42").Value = "10.1 °C"
$ws.Range("E43").Value = "2026-02-10 17:20:24"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "89%"
$ws.Range("H43").NumberFormat = "General"
$ws.Range("K43").Value = "9.1 MJ/m2"
$ws.Range("O43").Value = "8.9 °C"
$ws.Range("E44").Value = "2026-02-10 17:20:27"
$ws.Range("I44").Value = "18.4 mm"
$ws.Range("O44").Value = "0.2 °C"
$ws.Range("E45").Value = "2026-02-10 17:20:29"
$ws.Range("I45").Value = "24.4 mm"
$ws.Range("J45").Value = "1006.0 hPa"
$ws.Range("E46").Value = "2026-02-10 17:20:32"
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = "84%"
$ws.Range("H46").NumberFormat = "General"
$ws.Range("J46").Value = "1006.3 hPa"
$ws.Range("O46").Value = "13.8 °C"
